# "adding documents replacing mechanism"
# - Metryka!B7 title text changes from "Lista elementów" to "Spis dokumentacji wyrobu"
# - Lista_dokumentów column B is widened (was auto bestFit ~10 chars) to fit the new, longer header
# - cursor/selection position is moved on both sheets (cosmetic, matches the saved view state)

$wb = $excel.ActiveWorkbook

$wsMetryka = $wb.Worksheets.Item("Metryka")
$wsLista   = $wb.Worksheets.Item("Lista_dokumentów")

# Rename the document-title placeholder cell.
$wsMetryka.Range("B7").Value = "Spis dokumentacji wyrobu"

# Widen column B on the documents-list sheet so the new/longer strings fit.
$wsLista.Columns.Item(2).ColumnWidth = 21.5

# Restore the on-screen selections as last left by the editor.
$wsMetryka.Range("E12").Select()
$wsLista.Range("C5").Select()
